# Auto-generated edit script applying the Leve profit-data refresh
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 574.75
$ws.Range("I18").Value = 566.6667
$ws.Range("J18").Value = 599
$ws.Range("K18").Value = 566.6667
$ws.Range("L18").Value = 599
$ws.Range("M18").Value = -282.6667
$ws.Range("N18").Value = -1167

$ws.Range("H92").Value = 276.66666
$ws.Range("J92").Value = 531.5
$ws.Range("L92").Value = 531.5
$ws.Range("N92").Value = -3027.5

$ws.Range("H125").Value = 2399.5
$ws.Range("I125").Value = 2399.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 21595.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -19135.5
$ws.Range("N125").ClearContents()

$ws.Range("H127").Value = 796.5
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H131").Value = 7999.5
$ws.Range("I131").Value = 7999
$ws.Range("K131").Value = 23997
$ws.Range("M131").Value = -18957

$ws.Range("H138").Value = 2821.3455
$ws.Range("I138").Value = 2103.8333
$ws.Range("J138").Value = 3021.5813
$ws.Range("K138").Value = 6311.499899999999
$ws.Range("L138").Value = 9064.743899999999
$ws.Range("M138").Value = -1171.499899999999
$ws.Range("N138").Value = -19344.7439

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4197
$ws.Range("I105").Value = 7000
$ws.Range("J105").Value = 1394
$ws.Range("K105").Value = 7000
$ws.Range("L105").Value = 1394
$ws.Range("M105").Value = -5253
$ws.Range("N105").Value = -4888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2263.3572
$ws.Range("I31").Value = 2314.4167
$ws.Range("K31").Value = 2314.4167
$ws.Range("M31").Value = -2019.4167

$ws.Range("H34").Value = 2263.3572
$ws.Range("I34").Value = 2314.4167
$ws.Range("K34").Value = 2314.4167
$ws.Range("M34").Value = -2112.4167

$ws.Range("H107").Value = 1139.3529
$ws.Range("I107").Value = 519.25
$ws.Range("K107").Value = 519.25
$ws.Range("M107").Value = 1400.75

$ws.Range("H132").Value = 3556.8572
$ws.Range("I132").Value = 2950
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 8850
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -6320
$ws.Range("N132").Value = -16458.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 124
$ws.Range("J2").Value = 69.666664
$ws.Range("L2").Value = 417.999984
$ws.Range("N2").Value = -643.999984

$ws.Range("H137").Value = 2485.5715
$ws.Range("I137").Value = 2233.1667
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 6699.500100000001
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -1599.500100000001
$ws.Range("N137").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.85714
$ws.Range("I2").Value = 87.333336
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 87.333336
$ws.Range("L2").Value = 70
$ws.Range("M2").Value = 25.666664
$ws.Range("N2").Value = -296

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1630.3334
$ws.Range("J122").Value = 1495
$ws.Range("L122").Value = 4485
$ws.Range("N122").Value = -9385

$ws.Range("H132").Value = 2266.0527
$ws.Range("I132").Value = 1167.7142
$ws.Range("K132").Value = 3503.1426
$ws.Range("M132").Value = -973.1425999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2888
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 900
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 900
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1114

$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250

$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251

$ws.Range("H122").Value = 5500.375
$ws.Range("I122").Value = 4999.6665
$ws.Range("K122").Value = 14998.9995
$ws.Range("M122").Value = -12548.9995

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 4284
$ws.Range("I132").Value = 1904
$ws.Range("J132").Value = 5474
$ws.Range("K132").Value = 5712
$ws.Range("L132").Value = 16422
$ws.Range("M132").Value = -3182
$ws.Range("N132").Value = -21482

$ws.Range("H136").Value = 6330.5
$ws.Range("I136").Value = 7329.6665
$ws.Range("J136").Value = 5331.3335
$ws.Range("K136").Value = 21988.9995
$ws.Range("L136").Value = 15994.0005
$ws.Range("M136").Value = -19438.9995
$ws.Range("N136").Value = -21094.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 126697.75
$ws.Range("I62").Value = 167797
$ws.Range("K62").Value = 167797
$ws.Range("M62").Value = -167173

$ws.Range("H65").Value = 126697.75
$ws.Range("I65").Value = 167797
$ws.Range("K65").Value = 838985
$ws.Range("M65").Value = -835865

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H122").Value = 2082.7144
$ws.Range("I122").Value = 2096.5
$ws.Range("K122").Value = 6289.5
$ws.Range("M122").Value = -3839.5

$ws.Range("H126").Value = 1856.7273
$ws.Range("I126").Value = 1931.3
$ws.Range("J126").Value = 1111
$ws.Range("K126").Value = 5793.9
$ws.Range("L126").Value = 3333
$ws.Range("M126").Value = -3323.9
$ws.Range("N126").Value = -8273

$ws.Range("H132").Value = 3265.6
$ws.Range("I132").Value = 2716.6
$ws.Range("J132").Value = 3631.6
$ws.Range("K132").Value = 8149.799999999999
$ws.Range("L132").Value = 10894.8
$ws.Range("M132").Value = -5619.799999999999
$ws.Range("N132").Value = -15954.8

Write-Host "Applied Leve profit-data refresh across ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets"
